$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.402.73"
$ws.Range("E2").Value = "  +1.53%  "

$ws.Range("D3").Value = "3.143.39"
$ws.Range("E3").Value = "  +0.60%  "

$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.87"
$ws.Range("E5").Value = "  -1.00%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.16"
$ws.Range("E6").Value = "  -0.32%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("D8").Value = "3.137.04"
$ws.Range("E8").Value = "  +0.34%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.523"
$ws.Range("E9").Value = "  +0.82%  "

$ws.Range("E10").Value = "  -0.07%  "

$ws.Range("E11").Value = "  +2.76%  "

$ws.Range("E12").Value = "  +0.54%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000254"
$ws.Range("E13").Value = "  +1.54%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.48"
$ws.Range("E14").Value = "  +0.82%  "

$ws.Range("D15").Value = "3.662.43"
$ws.Range("E15").Value = "  +0.75%  "

$ws.Range("E16").Value = "  +2.81%  "

$ws.Range("D17").Value = "64.480.08"
$ws.Range("E17").Value = "  +1.48%  "

$ws.Range("D18").Value = "3.145.12"
$ws.Range("E18").Value = "  +0.65%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.88"
$ws.Range("E19").Value = "  +1.28%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "482.31"
$ws.Range("E20").Value = "  +1.44%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.64"
$ws.Range("E21").Value = "  +0.68%  "

$ws.Range("E23").Value = "  -0.97%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.03"
$ws.Range("E24").Value = "  +5.20%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.46"
$ws.Range("E25").Value = "  -0.06%  "

$ws.Range("E26").Value = "  +0.00%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.76"
$ws.Range("E27").Value = "  -0.70%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.43"
$ws.Range("E28").Value = "  +0.30%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.13"
$ws.Range("E29").Value = "  +4.41%  "

$ws.Range("E30").Value = "  +0.71%  "

$ws.Range("E31").Value = "  -4.93%  "

$ws.Range("E32").Value = "  +0.12%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.93"
$ws.Range("E33").Value = "  +2.94%  "

$ws.Range("E34").Value = "  +0.97%  "

$ws.Range("E35").Value = "  -1.61%  "

$ws.Range("E36").Value = "  +2.14%  "

$ws.Range("E37").Value = "  -0.24%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "52.90"
$ws.Range("E38").Value = "  +0.15%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.04"
$ws.Range("E39").Value = "  +3.83%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "441.55"
$ws.Range("E40").Value = "  -2.79%  "

$ws.Range("E41").Value = "  +1.20%  "

$ws.Range("E42").Value = "  +0.97%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.27"
$ws.Range("E43").Value = "  -0.46%  "

$ws.Range("D44").Value = "2.875.12"
$ws.Range("E44").Value = "  +1.27%  "

$ws.Range("E45").Value = "  -0.32%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.25"
$ws.Range("E46").Value = "  -1.21%  "

$ws.Range("E47").Value = "  +2.38%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.999"
$ws.Range("E48").Value = "  +0.01%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "26.06"
$ws.Range("E49").Value = "  +0.04%  "

$ws.Range("E50").Value = "  +0.62%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "121.60"
$ws.Range("E51").Value = "  +2.21%  "
